$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd phone number for the second student (a digit was missing).
$ws.Range("H3").Value = 777222333

# Leave the cursor on the cell that was just edited.
$null = $ws.Range("H3").Select()
